# "double header fixed in change password error page"
#
# 1) The footer/header "datetimeFigureOut" date placeholder text is
#    updated from 4/21/2011 to 4/25/2011 on the slide master and on
#    every slide layout (12 locations total).
# 2) Slide 10 ("Fun Facts") content placeholder text is corrected:
#      - "File count:  " + "27 " + "files" runs collapsed into a
#        single run "File count:  27 files"
#      - Comment lines:  60   -> 87
#      - Blank lines:  184    -> split into "Blank lines:  " + "185"
#      - CSS:  146            -> 147
#      - PHP:  1255           -> 1248
#      - Total:  1645         -> 1667

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            $rng = $shp.TextFrame.TextRange
            if ($rng.Text -eq "4/21/2011") {
                $rng.Text = "4/25/2011"
            }
        }
    }
}

# --- Update the date placeholder on the slide master ---
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# --- Update the date placeholder on every slide layout ---
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShapes $layout.Shapes
}

# --- Fix the "Fun Facts" slide (slide 10) statistics ---
$slide = $p.Slides.Item(10)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 1: "File count:  " / "27 " / "files" -> one run
$para1 = $tr.Paragraphs(1)
$sel1 = $tr.Characters($para1.Start, $para1.Length - 1)
$sel1.Text = "File count:  27 files"

# Paragraph 3: Comment lines 60 -> 87
[void]$tr.Paragraphs(3).Replace("60", "87", 1, 0, 0)

# Paragraph 4: Blank lines 184 -> 185, split into two runs
$para4 = $tr.Paragraphs(4)
$sel4 = $tr.Characters($para4.Start + 14, 3)
$sel4.Text = "185"

# Paragraph 5: CSS 146 -> 147
[void]$tr.Paragraphs(5).Replace("146", "147", 1, 0, 0)

# Paragraph 6: PHP 1255 -> 1248
[void]$tr.Paragraphs(6).Replace("1255", "1248", 1, 0, 0)

# Paragraph 7: Total 1645 -> 1667
[void]$tr.Paragraphs(7).Replace("1645", "1667", 1, 0, 0)
